$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2878.2632
$ws.Range("I116").Value = 2483.3845
$ws.Range("K116").Value = 2483.3845
$ws.Range("M116").Value = 958.6154999999999
$ws.Range("H132").Value = 18528332
$ws.Range("I132").Value = 22233120
$ws.Range("K132").Value = 66699360
$ws.Range("M132").Value = -66696830
$ws.Range("H134").Value = 59640
$ws.Range("J134").Value = 59640
$ws.Range("L134").Value = 59640
$ws.Range("N134").Value = -69780
$ws.Range("H136").Value = 78333.336
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 78333.336
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 78333.336
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -88533.336
$ws.Range("H137").Value = 1108.7273
$ws.Range("I137").Value = 985.3333
$ws.Range("K137").Value = 2955.9999
$ws.Range("M137").Value = -405.9998999999998
$ws.Range("H141").Value = 1591.5454
$ws.Range("I141").Value = 939.7143
$ws.Range("J141").Value = 2732.25
$ws.Range("K141").Value = 2819.1429
$ws.Range("L141").Value = 8196.75
$ws.Range("M141").Value = 2360.8571
$ws.Range("N141").Value = -18556.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 10006
$ws.Range("J9").Value = 10006
$ws.Range("L9").Value = 10006
$ws.Range("N9").Value = -10346
$ws.Range("H10").Value = 260
$ws.Range("J10").Value = 260
$ws.Range("L10").Value = 260
$ws.Range("N10").Value = -600
$ws.Range("H20").Value = 10006
$ws.Range("J20").Value = 10006
$ws.Range("L20").Value = 10006
$ws.Range("N20").Value = -10546
$ws.Range("H32").Value = 2930.9177
$ws.Range("I32").Value = 2622.3623
$ws.Range("K32").Value = 2622.3623
$ws.Range("M32").Value = -2335.3623
$ws.Range("H42").Value = 4277
$ws.Range("I42").Value = 1400
$ws.Range("K42").Value = 1400
$ws.Range("M42").Value = -914
$ws.Range("H45").Value = 1337.4706
$ws.Range("I45").Value = 1449
$ws.Range("K45").Value = 1449
$ws.Range("M45").Value = -1072
$ws.Range("H92").Value = 1268943.8
$ws.Range("J92").Value = 1268943.8
$ws.Range("L92").Value = 1268943.8
$ws.Range("N92").Value = -1273935.8
$ws.Range("H122").Value = 1958.1111
$ws.Range("I122").Value = 1958.1111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5874.3333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3424.3333
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1807.9688
$ws.Range("I132").Value = 1484.9131
$ws.Range("K132").Value = 4454.7393
$ws.Range("M132").Value = -1924.7393

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 50002004
$ws.Range("I94").Value = 83334670
$ws.Range("J94").Value = 3009.5
$ws.Range("K94").Value = 83334670
$ws.Range("L94").Value = 3009.5
$ws.Range("M94").Value = -83334219
$ws.Range("N94").Value = -3911.5
$ws.Range("H134").Value = 7730.684
$ws.Range("I134").Value = 1125.6
$ws.Range("K134").Value = 3376.8
$ws.Range("M134").Value = -841.7999999999997

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 6667168.5
$ws.Range("I12").Value = 752.5
$ws.Range("K12").Value = 752.5
$ws.Range("M12").Value = -582.5
$ws.Range("H31").Value = 1618.8667
$ws.Range("I31").Value = 1198.7646
$ws.Range("K31").Value = 1198.7646
$ws.Range("M31").Value = -903.7646
$ws.Range("H34").Value = 1618.8667
$ws.Range("I34").Value = 1198.7646
$ws.Range("K34").Value = 1198.7646
$ws.Range("M34").Value = -996.7646
$ws.Range("H132").Value = 6197.0386
$ws.Range("I132").Value = 8575.799999999999
$ws.Range("J132").Value = 2953.2727
$ws.Range("K132").Value = 25727.4
$ws.Range("L132").Value = 8859.8181
$ws.Range("M132").Value = -23197.4
$ws.Range("N132").Value = -13919.8181
$ws.Range("H134").Value = 2013.6786
$ws.Range("I134").Value = 1992.9546
$ws.Range("J134").Value = 2089.6667
$ws.Range("K134").Value = 5978.8638
$ws.Range("L134").Value = 6269.000100000001
$ws.Range("M134").Value = -3443.8638
$ws.Range("N134").Value = -11339.0001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1541.8
$ws.Range("I11").Value = 1854.75
$ws.Range("K11").Value = 5564.25
$ws.Range("M11").Value = -5424.25
$ws.Range("H105").Value = 11999.667
$ws.Range("J105").Value = 11999.667
$ws.Range("L105").Value = 35999.001
$ws.Range("N105").Value = -41241.001
$ws.Range("H116").Value = 2025
$ws.Range("I116").Value = 300
$ws.Range("J116").Value = 2600
$ws.Range("K116").Value = 900
$ws.Range("L116").Value = 7800
$ws.Range("M116").Value = 2542
$ws.Range("N116").Value = -14684
$ws.Range("H131").Value = 13334502
$ws.Range("I131").Value = 142857400
$ws.Range("J131").Value = 1262.4412
$ws.Range("K131").Value = 428572200
$ws.Range("L131").Value = 3787.3236
$ws.Range("M131").Value = -428567160
$ws.Range("N131").Value = -13867.3236

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4017
$ws.Range("I36").Value = 4017
$ws.Range("K36").Value = 4017
$ws.Range("M36").Value = -3532
$ws.Range("H43").Value = 5263.3335
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 6895
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 6895
$ws.Range("M43").Value = -1849
$ws.Range("N43").Value = -7197
$ws.Range("H46").Value = 18698.916
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 18698.916
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 18698.916
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -19010.916
$ws.Range("H102").Value = 2215.8647
$ws.Range("I102").Value = 2554.889
$ws.Range("J102").Value = 1894.6842
$ws.Range("K102").Value = 2554.889
$ws.Range("L102").Value = 1894.6842
$ws.Range("M102").Value = -932.8890000000001
$ws.Range("N102").Value = -5138.6842

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4000
$ws.Range("J5").Value = 4000
$ws.Range("L5").Value = 4000
$ws.Range("N5").Value = -4226
$ws.Range("H46").Value = 2750
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2750
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2750
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3126
$ws.Range("H61").Value = 1539.8422
$ws.Range("I61").Value = 1521.5
$ws.Range("J61").Value = 1571.2858
$ws.Range("K61").Value = 1521.5
$ws.Range("L61").Value = 1571.2858
$ws.Range("M61").Value = -1319.5
$ws.Range("N61").Value = -1975.2858
$ws.Range("H93").Value = 1083.2632
$ws.Range("I93").Value = 638.8
$ws.Range("J93").Value = 2750
$ws.Range("K93").Value = 638.8
$ws.Range("L93").Value = 2750
$ws.Range("M93").Value = 609.2
$ws.Range("N93").Value = -5246
$ws.Range("H113").Value = 1539.8422
$ws.Range("I113").Value = 1521.5
$ws.Range("J113").Value = 1571.2858
$ws.Range("K113").Value = 1521.5
$ws.Range("L113").Value = 1571.2858
$ws.Range("M113").Value = 648.5
$ws.Range("N113").Value = -5911.2858
$ws.Range("H122").Value = 7817053
$ws.Range("I122").Value = 10005315
$ws.Range("K122").Value = 30015945
$ws.Range("M122").Value = -30013495

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1000000000
$ws.Range("I26").Value = 1000000000
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1000000000
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -999999707
$ws.Range("N26").ClearContents()
$ws.Range("H62").Value = 125006750
$ws.Range("I62").Value = 166674340
$ws.Range("K62").Value = 166674340
$ws.Range("M62").Value = -166673716
$ws.Range("H65").Value = 125006750
$ws.Range("I65").Value = 166674340
$ws.Range("K65").Value = 833371700
$ws.Range("M65").Value = -833368580
$ws.Range("H107").Value = 398.8
$ws.Range("J107").Value = 250
$ws.Range("L107").Value = 750
$ws.Range("N107").Value = -4590
$ws.Range("H122").Value = 10001802
$ws.Range("I122").Value = 12382779
$ws.Range("K122").Value = 37148337
$ws.Range("M122").Value = -37145887
$ws.Range("H132").Value = 2054.6428
$ws.Range("I132").Value = 1721.7
$ws.Range("J132").Value = 2887
$ws.Range("K132").Value = 5165.1
$ws.Range("L132").Value = 8661
$ws.Range("M132").Value = -2635.1
$ws.Range("N132").Value = -13721
$ws.Range("H136").Value = 932.8889
$ws.Range("I136").Value = 716
$ws.Range("J136").Value = 1366.6666
$ws.Range("K136").Value = 2148
$ws.Range("L136").Value = 4099.9998
$ws.Range("M136").Value = 402
$ws.Range("N136").Value = -9199.9998
